$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("SoFCtMbCtPR ")

# Row 24 currently holds "hydrogen" -> rename it to "hydrogen combustion turbine"
$ws.Range("A24").Value = "hydrogen combustion turbine"

# Add a new row 25 for "hydrogen combined cycle" with the same share value as the rest
$ws.Range("A25").Value = "hydrogen combined cycle"
$ws.Range("B25").Value = 1.075

# Update the active sheet / selection to match the saved state (About sheet active,
# selection on the SoFCtMbCtPR sheet moved to B26)
$ws.Range("B26").Select()
$aboutWs = $wb.Worksheets.Item("About")
$aboutWs.Activate()
